# Update "想去人数" (interested count) figures in both the "展览" sheet
# and the "全部类型" sheet, which duplicate the same events.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1662
$ws1.Range("F3").Value = 866
$ws1.Range("F7").Value = 810
$ws1.Range("F8").Value = 831
$ws1.Range("F9").Value = 1539
$ws1.Range("F10").Value = 309
$ws1.Range("F14").Value = 205
$ws1.Range("F15").Value = 61
$ws1.Range("F16").Value = 516
$ws1.Range("F17").Value = 66
$ws1.Range("F19").Value = 11
$ws1.Range("F24").Value = 57
$ws1.Range("F28").Value = 196
$ws1.Range("F30").Value = 376

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1662
$ws4.Range("F5").Value = 866
$ws4.Range("F10").Value = 810
$ws4.Range("F11").Value = 831
$ws4.Range("F12").Value = 1539
$ws4.Range("F13").Value = 309
$ws4.Range("F17").Value = 205
$ws4.Range("F18").Value = 61
$ws4.Range("F19").Value = 516
$ws4.Range("F20").Value = 66
$ws4.Range("F23").Value = 11
$ws4.Range("F32").Value = 57
$ws4.Range("F37").Value = 196
$ws4.Range("F43").Value = 376

$wb.Save()
